# Append 14 new rows (317-330) of landscaping observation data to Sheet1,
# continuing the existing table for 2025-06-24 and 2025-06-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row data (excludes column A [date] and column F [formula], handled separately) ---
# Column order: B, C, D, E, G, H, I, J, K, L, M, N, O, P, Q, R, S, T
$rows = @(
    ,@(317, "Flowering", "Large", 74, 96, 0, 0.1, "No", 2, "Bright", 10, 0.52, 71, 30.16, 9, 0.14000000000000001, 9.9, 54, 7)
    ,@(318, "Nonflowering", "Medium", 74, 96, 0, 0.2, "No", 3, "Bright", 10, 0.52, 71, 30.16, 9, 0.14000000000000001, 9.9, 54, 7)
    ,@(319, "Nonflowering", "Small", 74, 96, 0, 0, "No", 3, "Neutral", 10, 0.52, 71, 30.16, 9, 0.14000000000000001, 9.9, 54, 7)
    ,@(320, "Nonflowering", "Medium", 74, 96, 0, 0, "No", 3, "Bright", 10, 0.52, 71, 30.16, 9, 0.14000000000000001, 9.9, 54, 7)
    ,@(321, "Nonflowering", "Medium", 74, 96, 0, 0.1, "No", 3, "Bright", 10, 0.52, 71, 30.16, 9, 0.14000000000000001, 9.9, 54, 7)
    ,@(322, "Nonflowering", "Large", 74, 96, 0, 0.5, "No", 4, "Bright", 10, 0.52, 71, 30.16, 9, 0.14000000000000001, 9.9, 54, 7)
    ,@(323, "Tree", "Medium", 74, 96, 0, 0.3, "No", 1, "Bright", 10, 0.52, 71, 30.16, 9, 0.14000000000000001, 9.9, 54, 7)
    ,@(324, "Flowering", "Large", 73, 91, 0, 0, "No", 2, "Neutral", 8, 0.5, 70, 30.14, 14, 0.61, 9.9, 62, 6)
    ,@(325, "Nonflowering", "Medium", 73, 91, 0, 0, "No", 3, "Neutral", 8, 0.5, 70, 30.14, 14, 0.61, 9.9, 62, 6)
    ,@(326, "Nonflowering", "Small", 73, 91, 0, 0, "No", 3, "Bright", 8, 0.5, 70, 30.14, 14, 0.61, 9.9, 62, 6)
    ,@(327, "Nonflowering", "Medium", 73, 91, 0, 0.1, "No", 3, "Bright", 8, 0.5, 70, 30.14, 14, 0.61, 9.9, 62, 6)
    ,@(328, "Nonflowering", "Medium", 73, 91, 0, 0.1, "No", 3, "Bright", 8, 0.5, 70, 30.14, 14, 0.61, 9.9, 62, 6)
    ,@(329, "Nonflowering", "Large", 73, 91, 0, 0.4, "No", 4, "Neutral", 8, 0.5, 70, 30.14, 14, 0.61, 9.9, 62, 6)
    ,@(330, "Tree", "Medium", 73, 91, 0, 0.2, "No", 1, "Bright", 8, 0.5, 70, 30.14, 14, 0.61, 9.9, 62, 6)
)

$dates = @{317=45832; 318=45832; 319=45832; 320=45832; 321=45832; 322=45832; 323=45832; 324=45833; 325=45833; 326=45833; 327=45833; 328=45833; 329=45833; 330=45833}

# --- Column A: copy the date cell format (style) down from the last existing row, then fill in values ---
$ws.Range("A316").Copy()
$ws.Range("A317:A330").PasteSpecial(-4122)
$excel.CutCopyMode = $false

foreach ($rownum in 317..330) {
    $ws.Cells.Item($rownum, 1).Value = $dates[$rownum]
}

# --- Columns B-E, G-T ---
$colMap = @(2, 3, 4, 5, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)

foreach ($r in $rows) {
    $rownum = $r[0]
    for ($i = 0; $i -lt $colMap.Length; $i++) {
        $col = $colMap[$i]
        $val = $r[$i + 1]
        $ws.Cells.Item($rownum, $col).Value = $val
    }
}

# --- Column F: ABS(D-E) temperature-difference formula, matching the existing pattern ---
foreach ($rownum in 317..330) {
    $ws.Cells.Item($rownum, 6).Formula = "=ABS(D" + $rownum + "-E" + $rownum + ")"
}

# --- View state: scroll window down to the newly added rows and select N324:N330 ---
$ws.Activate()
$ws.Range("N324:N330").Select()
$excel.ActiveWindow.ScrollRow = 306
$excel.ActiveWindow.ScrollColumn = 1
